$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startSerial = 44313
for ($i = 0; $i -lt 6; $i++) {
    $row = 239 + $i
    $serial = $startSerial + $i
    $ws.Cells.Item(238, 1).Copy($ws.Cells.Item($row, 1))
    $ws.Cells.Item($row, 1).Value = $serial
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0
}
